$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 172 ("basal diameter" / QUALITY / 1) entirely, shifting all
# subsequent rows up by one.
$ws.Rows(172).Delete()
